# LF Denoising results (noisy)
# Fills in the previously-empty "PB-LAB402" (column I) results on sheet x2,
# and the previously-empty column-E results for rows 15/16 on sheet x4.
# Downstream AVERAGE() formulas (I17/I18/I20 on x2, E19/E20 on x4) recompute
# automatically, turning their prior #DIV/0! errors into real numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "x3" - only the selection moves (no data changes on this sheet)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("x3")
$ws3.Range("E30").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "x4" - fill in E15/E16 (noisy-LF results), matching the number
# format used elsewhere in the column (0.000, centered, borderless - the
# new style distinct from the bordered neighbours).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("x4")

$ws4.Range("E15").NumberFormat = "0.000"
$ws4.Range("E15").HorizontalAlignment = -4108
$ws4.Range("E15").VerticalAlignment = -4108
$ws4.Range("E15").Borders.LineStyle = 0
$ws4.Range("E15").Value = 29.345099999999999

$ws4.Range("E16").NumberFormat = "0.000"
$ws4.Range("E16").HorizontalAlignment = -4108
$ws4.Range("E16").VerticalAlignment = -4108
$ws4.Range("E16").Borders.LineStyle = 0
$ws4.Range("E16").Value = 26.133199999999999

$ws4.Range("G30").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "x2" - fill in I2:I8 (noisy-LF PB-LAB402 results); this sheet
# becomes the active tab with I8 selected, as in the target workbook.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("x2")

$ws2.Range("I2").Value = 31.9481
$ws2.Range("I3").Value = 26.151499999999999
$ws2.Range("I4").Value = 28.038900000000002
$ws2.Range("I5").Value = 36.229599999999998
$ws2.Range("I6").Value = 31.1417
$ws2.Range("I7").Value = 33.225299999999997
$ws2.Range("I8").Value = 27.573

$ws2.Range("I8").Select() | Out-Null
